# Update column G ("K" - strikeouts) with newly computed values (s_vals)
# for each game row, per commit: "regen save_data to use K instead of Strike#,
# regen std/mean, calc and write s_vals"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 2
    3  = 5
    4  = 3
    5  = 4
    6  = 6
    7  = 7
    8  = 2
    9  = 2
    10 = 2
    11 = 8
    12 = 6
    13 = 2
    14 = 3
    15 = 3
    16 = 8
    17 = 5
    18 = 5
    19 = 7
    20 = 6
    21 = 4
    22 = 5
    23 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
